$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B20").Value = 2
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("I20").Value = 2
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 1
$ws.Range("N20").Value = 1
$ws.Range("O20").Value = 2
$ws.Range("P20").Value = 1
$ws.Range("Q20").Value = 2
$ws.Range("R20").Value = 2
$ws.Range("S20").Value = 1
$ws.Range("T20").Value = 1
$ws.Range("U20").Value = 2
$ws.Range("V20").Value = 1
$ws.Range("W20").Value = 1
$ws.Range("X20").Value = 1
$ws.Range("Y20").Value = 2
$ws.Range("Z20").Value = 2
$ws.Range("AA20").Value = 3
$ws.Range("AB20").Value = 1
$ws.Range("AC20").Value = 1
$ws.Range("AD20").Value = 2
$ws.Range("AE20").Value = 1
$ws.Range("AF20").Value = 2
$ws.Range("AG20").Value = 2
$ws.Range("AH20").Value = 1
$ws.Range("AI20").Value = 4
$ws.Range("AK20").Value = 3
$ws.Range("A28").Value = 'Consejo Nacional de Investigaciones CientÃ­ficas y TÃ©cnicas'
$ws.Range("B39").Value = 3
$ws.Range("C39").Value = 4
$ws.Range("D39").Value = 3
$ws.Range("E39").Value = 1
$ws.Range("F39").Value = 1
$ws.Range("G39").Value = 1
$ws.Range("H39").Value = 2
$ws.Range("I39").Value = 3
$ws.Range("J39").Value = 1
$ws.Range("L39").Value = 1
$ws.Range("N39").Value = 1
$ws.Range("O39").Value = 2
$ws.Range("R39").Value = 1
$ws.Range("S39").Value = 2
$ws.Range("T39").Value = 1
$ws.Range("U39").Value = 2
$ws.Range("V39").Value = 1
$ws.Range("X39").Value = 2
$ws.Range("Y39").Value = 2
$ws.Range("Z39").Value = 2
$ws.Range("AA39").Value = 3
$ws.Range("AB39").Value = 2
$ws.Range("AD39").Value = 3
$ws.Range("AE39").Value = 1
$ws.Range("AF39").Value = 2
$ws.Range("AG39").Value = 2
$ws.Range("AI39").Value = 5
$ws.Range("AK39").Value = 3
$ws.Range("B72").Value = 1
$ws.Range("C72").Value = 1
$ws.Range("D72").Value = 3
$ws.Range("H72").Value = 1
$ws.Range("I72").Value = 1
$ws.Range("K72").Value = 1
$ws.Range("L72").Value = 1
$ws.Range("N72").Value = 1
$ws.Range("P72").Value = 1
$ws.Range("Q72").Value = 1
$ws.Range("W72").Value = 1
$ws.Range("X72").Value = 1
$ws.Range("Z72").Value = 2
$ws.Range("AA72").Value = 2
$ws.Range("AB72").Value = 1
$ws.Range("AC72").Value = 1
$ws.Range("AD72").Value = 1
$ws.Range("AF72").Value = 2
$ws.Range("AG72").Value = 3
$ws.Range("AH72").Value = 2
$ws.Range("AI72").Value = 2
$ws.Range("AK72").Value = 2
$ws.Range("B90").Value = 3
$ws.Range("C90").Value = 2
$ws.Range("D90").Value = 1
$ws.Range("E90").Value = 0
$ws.Range("F90").Value = 0
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 2
$ws.Range("I90").Value = 2
$ws.Range("N90").Value = 0
$ws.Range("O90").Value = 1
$ws.Range("R90").Value = 2
$ws.Range("U90").Value = 1
$ws.Range("X90").Value = 1
$ws.Range("Y90").Value = 1
$ws.Range("AA90").Value = 3
$ws.Range("AD90").Value = 2
$ws.Range("AE90").Value = 1
$ws.Range("AF90").Value = 1
$ws.Range("AG90").Value = 0
$ws.Range("AI90").Value = 3
$ws.Range("AK90").Value = 3
